# Applies the "Updated symbol list" crypto-price refresh (Sat Dec 31 21:46:44 UTC 2022).
# Each target cell is price/volume text (stored as inline strings, General format),
# so values are written with a leading apostrophe to force text-typed input, then
# the cell Style is reset to "Normal" so no stray number-format/quote-prefix style sticks.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'247.33"
$c.Style = "Normal"

$c = $ws.Range("E2")
$c.Value = "'0.87%"
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.Value = "'26.26"
$c.Style = "Normal"

$c = $ws.Range("E3")
$c.Value = "'4.53%"
$c.Style = "Normal"

$c = $ws.Range("D4")
$c.Value = "'5.089"
$c.Style = "Normal"

$c = $ws.Range("E4")
$c.Value = "'1.18%"
$c.Style = "Normal"

$c = $ws.Range("E5")
$c.Value = "'-0.33%"
$c.Style = "Normal"

$c = $ws.Range("D6")
$c.Value = "'6.479"
$c.Style = "Normal"

$c = $ws.Range("D7")
$c.Value = "'0.8128"
$c.Style = "Normal"

$c = $ws.Range("E7")
$c.Value = "'-0.11%"
$c.Style = "Normal"

$c = $ws.Range("D8")
$c.Value = "'0.8447"
$c.Style = "Normal"

$c = $ws.Range("E8")
$c.Value = "'1.22%"
$c.Style = "Normal"

$c = $ws.Range("D9")
$c.Value = "'0.06991"
$c.Style = "Normal"

$c = $ws.Range("E9")
$c.Value = "'0.60%"
$c.Style = "Normal"

$c = $ws.Range("D10")
$c.Value = "'0.02807"
$c.Style = "Normal"

$c = $ws.Range("E10")
$c.Value = "'-1.15%"
$c.Style = "Normal"

$c = $ws.Range("D11")
$c.Value = "'0.09385"
$c.Style = "Normal"

$c = $ws.Range("E11")
$c.Value = "'-0.18%"
$c.Style = "Normal"

$c = $ws.Range("D12")
$c.Value = "'0.001511"
$c.Style = "Normal"

$c = $ws.Range("E12")
$c.Value = "'-0.65%"
$c.Style = "Normal"

$c = $ws.Range("D13")
$c.Value = "'0.0005960"
$c.Style = "Normal"

$c = $ws.Range("E13")
$c.Value = "'-0.23%"
$c.Style = "Normal"

$c = $ws.Range("D14")
$c.Value = "'0.006165"
$c.Style = "Normal"

$c = $ws.Range("E14")
$c.Value = "'0.82%"
$c.Style = "Normal"

$c = $ws.Range("D15")
$c.Value = "'3.606"
$c.Style = "Normal"

$c = $ws.Range("E15")
$c.Value = "'3.05%"
$c.Style = "Normal"

$c = $ws.Range("E16")
$c.Value = "'0.29%"
$c.Style = "Normal"

$c = $ws.Range("D18")
$c.Value = "'0.3113"
$c.Style = "Normal"

$c = $ws.Range("E18")
$c.Value = "'-2.24%"
$c.Style = "Normal"

$c = $ws.Range("D19")
$c.Value = "'0.1340"
$c.Style = "Normal"

$c = $ws.Range("E19")
$c.Value = "'0.23%"
$c.Style = "Normal"

$c = $ws.Range("D20")
$c.Value = "'0.03193"
$c.Style = "Normal"

$c = $ws.Range("E20")
$c.Value = "'-2.15%"
$c.Style = "Normal"

$c = $ws.Range("E21")
$c.Value = "'-1.28%"
$c.Style = "Normal"

$c = $ws.Range("D22")
$c.Value = "'3.759"
$c.Style = "Normal"

$c = $ws.Range("E22")
$c.Value = "'0.59%"
$c.Style = "Normal"

$c = $ws.Range("D23")
$c.Value = "'0.04647"
$c.Style = "Normal"

$c = $ws.Range("E23")
$c.Value = "'-0.50%"
$c.Style = "Normal"

$c = $ws.Range("E24")
$c.Value = "'-1.35%"
$c.Style = "Normal"

$c = $ws.Range("D25")
$c.Value = "'0.001249"
$c.Style = "Normal"

$c = $ws.Range("E25")
$c.Value = "'0.56%"
$c.Style = "Normal"

$c = $ws.Range("D26")
$c.Value = "'0.004576"
$c.Style = "Normal"

$c = $ws.Range("E26")
$c.Value = "'1.02%"
$c.Style = "Normal"

$c = $ws.Range("E27")
$c.Value = "'-0.94%"
$c.Style = "Normal"

$c = $ws.Range("E28")
$c.Value = "'0.00%"
$c.Style = "Normal"

$c = $ws.Range("D40")
$c.Value = "'0.03662"
$c.Style = "Normal"

$c = $ws.Range("E40")
$c.Value = "'-0.03%"
$c.Style = "Normal"

$c = $ws.Range("B41")
$c.Value = "'BKEXToken"
$c.Style = "Normal"

$c = $ws.Range("C41")
$c.Value = "'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$c.Style = "Normal"

$c = $ws.Range("D41")
$c.Value = "'0.1356"
$c.Style = "Normal"

$c = $ws.Range("E41")
$c.Value = "'-0.03%"
$c.Style = "Normal"

$c = $ws.Range("B42")
$c.Value = "'CEJI"
$c.Style = "Normal"

$c = $ws.Range("C42")
$c.Value = "'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$c.Style = "Normal"

$c = $ws.Range("D42")
$c.Value = "'0.002660"
$c.Style = "Normal"

$c = $ws.Range("E42")
$c.Value = "'-2.65%"
$c.Style = "Normal"

$c = $ws.Range("B43")
$c.Value = "'KickToken"
$c.Style = "Normal"

$c = $ws.Range("C43")
$c.Value = "'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$c.Style = "Normal"

$c = $ws.Range("D43")
$c.Value = "'0.003437"
$c.Style = "Normal"

$c = $ws.Range("E43")
$c.Value = "'-44.78%"
$c.Style = "Normal"

$c = $ws.Range("D44")
$c.Value = "'0.008276"
$c.Style = "Normal"

$c = $ws.Range("E44")
$c.Value = "'1.38%"
$c.Style = "Normal"

$c = $ws.Range("D45")
$c.Value = "'0.00005388"
$c.Style = "Normal"

$c = $ws.Range("E45")
$c.Value = "'1.85%"
$c.Style = "Normal"

$c = $ws.Range("E46")
$c.Value = "'0.10%"
$c.Style = "Normal"

$c = $ws.Range("E47")
$c.Value = "'-35.78%"
$c.Style = "Normal"

$c = $ws.Range("D48")
$c.Value = "'0.002432"
$c.Style = "Normal"

$c = $ws.Range("E48")
$c.Value = "'20.32%"
$c.Style = "Normal"

$c = $ws.Range("E49")
$c.Value = "'0.10%"
$c.Style = "Normal"

$c = $ws.Range("E50")
$c.Value = "'0.10%"
$c.Style = "Normal"
